# Add the latest quotations row (2025-10-24) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 50

# Column A: date value (serial 45954 = 2025-10-24).
# Set the value first, then copy the number format from the row above
# so the new cell keeps the same date-like numeric formatting/style.
$ws.Cells.Item($newRow, 1).Value = 45954
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

# Columns B-E: quotation values (kept as text with comma decimal separator,
# matching the existing data in the sheet).
$ws.Cells.Item($newRow, 2).Value = "21,7048"
$ws.Cells.Item($newRow, 3).Value = "15,5758"
$ws.Cells.Item($newRow, 4).Value = "15,5156"
$ws.Cells.Item($newRow, 5).Value = "15,5156"
